# Apply the updated cryptocurrency price/volume snapshot to Sheet1.
# Cells whose new text would otherwise be auto-coerced to a number by the
# COM Value setter are written via NumberFormat "@" (text) and then have their
# formatting cleared again so no stray style survives on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.079.47"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "2.685.30"
$ws.Range("E3").Value = "  -2.69%  "
$ws.Range("E4").Value = "  -0.27%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "551.48"
$c.ClearFormats()
$ws.Range("E5").Value = "  -4.46%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "157.81"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.46%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E7").Value = "  -0.10%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.590"
$c.ClearFormats()
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("E9").Value = "  -4.49%  "
$ws.Range("E10").Value = "  -2.81%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.369"
$c.ClearFormats()
$ws.Range("E11").Value = "  -4.66%  "
$ws.Range("E12").Value = "  -12.01%  "
$ws.Range("D13").Value = "3.159.22"
$ws.Range("E13").Value = "  -2.87%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.16"
$c.ClearFormats()
$ws.Range("E14").Value = "  -4.17%  "
$ws.Range("D15").Value = "62.944.13"
$ws.Range("E15").Value = "  -1.93%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000148"
$c.ClearFormats()
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("D17").Value = "2.684.19"
$ws.Range("E17").Value = "  -3.16%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "11.96"
$c.ClearFormats()
$ws.Range("E18").Value = "  -1.73%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.59"
$c.ClearFormats()
$ws.Range("E19").Value = "  -5.50%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "343.76"
$c.ClearFormats()
$ws.Range("E20").Value = "  -4.12%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.32"
$c.ClearFormats()
$ws.Range("E21").Value = "  -5.29%  "
$ws.Range("E22").Value = "  +0.09%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.505"
$c.ClearFormats()
$ws.Range("E23").Value = "  -4.55%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "63.61"
$c.ClearFormats()
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("E26").Value = "  +0.10%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.14"
$c.ClearFormats()
$ws.Range("E27").Value = "  -5.53%  "
$ws.Range("D28").Value = "0.0₃0860"
$ws.Range("E28").Value = "  -7.09%  "
$ws.Range("E29").Value = "  -2.25%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.34"
$c.ClearFormats()
$ws.Range("E30").Value = "  -2.58%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.04"
$c.ClearFormats()
$ws.Range("E31").Value = "  -4.48%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "166.34"
$c.ClearFormats()
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("E33").Value = "  -0.03%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.83"
$c.ClearFormats()
$ws.Range("E34").Value = "  -3.42%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "19.55"
$c.ClearFormats()
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("E36").Value = "  -5.76%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.78"
$c.ClearFormats()
$ws.Range("E37").Value = "  -3.42%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "339.12"
$c.ClearFormats()
$ws.Range("E38").Value = "  -3.58%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "6.22"
$c.ClearFormats()
$ws.Range("E39").Value = "  -3.04%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.935"
$c.ClearFormats()
$ws.Range("E40").Value = "  -7.06%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.95"
$c.ClearFormats()
$ws.Range("E41").Value = "  -5.65%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "38.24"
$c.ClearFormats()
$ws.Range("E42").Value = "  -2.10%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "20.43"
$c.ClearFormats()
$ws.Range("E43").Value = "  -5.68%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "20.83"
$c.ClearFormats()
$ws.Range("E44").Value = "  -7.58%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.618"
$c.ClearFormats()
$ws.Range("E45").Value = "  -2.05%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0562"
$c.ClearFormats()
$ws.Range("E46").Value = "  -5.30%  "
$ws.Range("E47").Value = "  -0.17%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "11.05"
$c.ClearFormats()
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  -3.99%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "129.26"
$c.ClearFormats()
$ws.Range("E50").Value = "  -5.45%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0241"
$c.ClearFormats()
$ws.Range("E51").Value = "  -4.99%  "
